# Updates the "Förändrad" (Changed) date column (C) for all data rows
# from serial date 45783 (2025-05-06) to 45784 (2025-05-07).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 2; $row -le 43; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45783) {
        $cell.Value = 45784
    }
}
